# Final-capstone-story-project.pptx edit
# Commit: "Update final project files"
#
# Slide 4 ("INTRODUCTION") has a "Content Placeholder 2" text box whose body
# contained a stray placeholder paragraph "Point3" right before the first
# real bullet ("What is the most used programming language today?").
# Remove that stray paragraph entirely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Locate the shape named "Content Placeholder 2" (defensive lookup by name,
# falls back to the known index if Name isn't resolvable for some reason).
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "Content Placeholder 2") {
        $shape = $s.Shapes.Item($i)
        break
    }
}
if ($shape -eq $null) {
    $shape = $s.Shapes.Item(3)
}

$tr = $shape.TextFrame.TextRange

# Find and delete the paragraph whose text is exactly "Point3".
# Note: TextRange.Text includes the trailing paragraph-mark ("`r"), so
# trim it before comparing.
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    $txt = $para.Text.TrimEnd("`r")
    if ($txt -eq "Point3") {
        $para.Delete()
        break
    }
}
